{"js": "// The only substantive textual change in this diff is in the \"Designed an\n// interchangeable controller\" bullet: \"Designed an interchangeable controller \"\n// becomes \"Designed interchangeable controllers \" (drop \"an \", pluralize\n// \"controller\" -> \"controllers\"). Every other hunk in the diff only adds\n// <w:proofErr> spell-check markers around existing words (Automation,\n// Ignition, Matlab, Flexlink, SureKap, AutoBagger, Jaclean) by splitting an\n// existing run into several runs with identical formatting \u2014 the visible\n// text and formatting are unchanged, so nothing needs to happen there.\n\nconst body = context.document.body;\n\nconst results = body.search(\"Designed an interchangeable controller \", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find target phrase to replace.\");\n}\n\nresults.items[0].insertText(\"Designed interchangeable controllers \", \"Replace\");\nawait context.sync();\n", "ps1": "# The only substantive textual change in this diff is in the \"Designed an\n# interchangeable controller\" bullet: \"Designed an interchangeable controller \"\n# becomes \"Designed interchangeable controllers \" (drop \"an \", pluralize\n# \"controller\" -> \"controllers\"). Every other hunk in the diff only adds\n# <w:proofErr> spell-check markers around existing words (Automation,\n# Ignition, Matlab, Flexlink, SureKap, AutoBagger, Jaclean) by splitting an\n# existing run into several runs with identical formatting - the visible\n# text and formatting are unchanged, so nothing needs to happen there.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"Designed an interchangeable controller \"\n$find.Replacement.Text = \"Designed interchangeable controllers \"\n$result = $find.Execute(\n    $null,   # FindText\n    $true,   # MatchCase\n    $false,  # MatchWholeWord\n    $false,  # MatchWildcards\n    $false,  # MatchSoundsLike\n    $false,  # MatchAllWordForms\n    $true,   # Forward\n    1,       # Wrap (wdFindContinue)\n    $false,  # Format\n    $null,   # ReplaceWith\n    2        # Replace (wdReplaceOne)\n)\n\nif (-not $result) {\n    throw \"Could not find target phrase to replace.\"\n}\n"}
